# Update the lattice-multiplication exercise table: every cell keeps its
# existing layout/formatting (single run, sz=32, 5 lines separated by
# <w:br/>) but gets a new "A x B" problem. The 5 lines of each cell follow
# a fixed template derived from the two 2-digit factors A and B:
#   line1: "A x B"
#   line2: "  B0    B1"      (tens digit, ones digit of B)
#   line3: "  ----"          (unchanged)
#   line4: "A0|    |"        (tens digit of A)
#   line5: "A1|    |"        (ones digit of A)
#
# NOTE: PowerShell's "+" operator performs numeric addition when both
# operands parse as numbers (e.g. "4" + "3" -> 7, NOT "43"), so all
# string building below uses string interpolation ("$a$b") rather than
# "+" concatenation to avoid accidental arithmetic.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# New problems, row-major (row1: col1..col3, row2: col1..col3, ...)
$problems = @(
    "23 x 43", "82 x 28", "45 x 71",
    "70 x 80", "48 x 28", "92 x 41",
    "68 x 25", "58 x 40", "35 x 37",
    "67 x 97", "88 x 71", "63 x 15",
    "35 x 39", "59 x 81", "85 x 17"
)

$vt = [char]11
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $problem = $problems[$idx]
        $idx = $idx + 1

        $parts = $problem.Split("x")
        $a = $parts[0].Trim()
        $b = $parts[1].Trim()

        $a0 = $a.Substring(0,1)
        $a1 = $a.Substring(1,1)
        $b0 = $b.Substring(0,1)
        $b1 = $b.Substring(1,1)

        $line1 = $problem
        $line2 = "  $b0    $b1"
        $line3 = "  ----"
        $line4 = "$a0|    |"
        $line5 = "$a1|    |"

        $full = "$line1$vt$line2$vt$line3$vt$line4$vt$line5"

        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        # Trim off the trailing end-of-cell marker (2 chars) so we don't
        # clobber the cell/row structure, then overwrite the text in place
        # (this preserves the existing run formatting, e.g. sz=32).
        $rng.MoveEnd(1, -2) | Out-Null
        $rng.Text = $full
    }
}
